# "Latest selenium version added, removed webdriver manager"
# The RunSheet's Browser column (B2:B14) previously mixed Chrome/Firefox/Edge
# run rows; now every run row drives Edge (Selenium Manager auto-resolves the
# Edge driver, so the separate WebDriverManager per-browser bookkeeping goes
# away). Set every data row in column B to "Edge".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunSheet")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "Edge"
}

# Mirror the author's last selection/cursor position left in the saved file.
$ws.Range("E9").Select() | Out-Null
